# The data rows (2-10) were reshuffled: each row keeps its identifying
# columns (A,B,C,E,F,G,H,I,J,K,L) but the varying measurement columns
# (D,M,N,O,P,Q,R,S,T) get swapped around between rows, per this mapping
# of destination row -> source row (row 8 is untouched):
#   2<-5, 3<-9, 4<-10, 5<-2, 6<-3, 7<-4, 9<-6, 10<-7

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Snapshot of the "varying" columns for every affected source row, taken
# BEFORE any writes happen (so overwriting one row doesn't clobber the
# value another row still needs to read).
$cols = @("D", "M", "N", "O", "P", "Q", "R", "S", "T")
$sourceRows = @(2, 3, 4, 5, 6, 7, 9, 10)

$snapshot = @{}
foreach ($r in $sourceRows) {
    $rowData = @{}
    foreach ($c in $cols) {
        $rowData[$c] = $ws.Range("$c$r").Value2
    }
    $snapshot[$r] = $rowData
}

# destination row -> source row
$mapping = @{
    2  = 5
    3  = 9
    4  = 10
    5  = 2
    6  = 3
    7  = 4
    9  = 6
    10 = 7
}

foreach ($dest in $mapping.Keys) {
    $src = $mapping[$dest]
    $rowData = $snapshot[$src]
    foreach ($c in $cols) {
        $ws.Range("$c$dest").Value = $rowData[$c]
    }
}
